# Adicionado a classe Item doacao e ajustes finais
# Remove the second user row (old row 2) and append a new user row at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old row 2 (user "dwad"); remaining rows shift up automatically.
$ws.Rows.Item(2).Delete()

# Append the new user's data as the new last row (row 5). A leading
# apostrophe forces these numeric-looking CPF/phone-style values to be
# stored as text (matching the existing column data), then ClearFormats
# drops the resulting quote-prefix cell style so the cells keep the sheet's
# default style.
$newRow5 = $ws.Range("A5:E5")
$ws.Range("A5").Value = "'321312"
$ws.Range("B5").Value = "'123213213"
$ws.Range("C5").Value = "'32131"
$ws.Range("D5").Value = "'2321"
$ws.Range("E5").Value = "'321"
$newRow5.ClearFormats()
